$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply the refreshed coin-list snapshot (prices / 1h-volume deltas) cell by
# cell. A handful of Price cells would otherwise be auto-coerced from text to
# a Number by Excel's smart entry (e.g. "0.999", "23.18"), unlike sibling
# cells such as "27.792.58" whose two dot separators keep them text naturally.
# Force the Text number format first on exactly those cells so every Price
# cell keeps matching the source sheet's string-typed storage.

$ws.Range("D2").Value = "27.792.58"
$ws.Range("D3").Value = "1.617.40"
$ws.Range("E3").Value = "  -1.46%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.12%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "209.37"
$ws.Range("E5").Value = "  -1.73%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.521"
$ws.Range("E6").Value = "  -0.69%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.998"
$ws.Range("E7").Value = "  -0.15%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "23.18"
$ws.Range("E8").Value = "  -1.86%  "
$ws.Range("E9").Value = "  -1.33%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0609"
$ws.Range("E10").Value = "  -1.00%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0878"
$ws.Range("E11").Value = "  -0.56%  "
$ws.Range("D12").Value = "1.845.79"
$ws.Range("E12").Value = "  -1.50%  "
$ws.Range("D13").Value = "1.608.75"
$ws.Range("E13").Value = "  -1.96%  "
$ws.Range("E14").Value = "  -2.21%  "
$ws.Range("E15").Value = "  -3.04%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.93"
$ws.Range("E16").Value = "  -1.48%  "
$ws.Range("D17").Value = "27.760.62"
$ws.Range("E17").Value = "  -0.83%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "228.28"
$ws.Range("E18").Value = "  -3.08%  "
$ws.Range("E19").Value = "  -0.92%  "
$ws.Range("E20").Value = "  -1.04%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.998"
$ws.Range("E21").Value = "  -0.19%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.29"
$ws.Range("E22").Value = "  -2.12%  "
$ws.Range("E23").Value = "  -5.74%  "
$ws.Range("E24").Value = "  -3.43%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "154.39"
$ws.Range("E25").Value = "  +2.15%  "
$ws.Range("E26").Value = "  -1.57%  "
$ws.Range("E27").Value = "  -0.51%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.43"
$ws.Range("E28").Value = "  -1.69%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.999"
$ws.Range("E29").Value = "  -0.17%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.17"
$ws.Range("E30").Value = "  -1.98%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0479"
$ws.Range("E31").Value = "  -1.11%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.39"
$ws.Range("E33").Value = "  -2.07%  "
$ws.Range("D34").Value = "1.386.30"
$ws.Range("E34").Value = "  -2.31%  "
$ws.Range("E35").Value = "  -1.16%  "
$ws.Range("E36").Value = "  +9.15%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.33"
$ws.Range("E37").Value = "  -1.14%  "
$ws.Range("E38").Value = "  -0.12%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.553"
$ws.Range("E39").Value = "  -0.92%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.847"
$ws.Range("E40").Value = "  -4.06%  "
$ws.Range("E41").Value = "  -1.20%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.997"
$ws.Range("E42").Value = "  -0.25%  "
$ws.Range("E43").Value = "  -2.79%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "65.41"
$ws.Range("E44").Value = "  -1.83%  "
$ws.Range("E45").Value = "  -1.32%  "
$ws.Range("B46").Value = "RocketPoolETH"
$ws.Range("C46").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D46").Value = "1.756.10"
$ws.Range("E46").Value = "  -1.51%  "
$ws.Range("B47").Value = "MXToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.15"
$ws.Range("E47").Value = "  -2.47%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "87.51"
$ws.Range("E48").Value = "  -0.41%  "
$ws.Range("B49").Value = "BabyDogeCoin"
$ws.Range("C49").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D49").Value = "0.0₆0103"
$ws.Range("E49").Value = "  -2.50%  "
$ws.Range("B50").Value = "Algorand"
$ws.Range("C50").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.100"
$ws.Range("E50").Value = "  -0.06%  "
$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0502"
$ws.Range("E51").Value = "  -0.77%  "
